# Adds a new slide ("Domain Controller") as slide 6, mirroring the layout/
# formatting conventions already used by slide 5 in this deck.

function EmuToPt($emu) {
    return $emu / 12700.0
}

function BgrFromRrggbb($rrggbb) {
    # PowerPoint's Font.Color.RGB (and similar .RGB longs) are ordered
    # 0x00BBGGRR, i.e. the reverse byte order of the usual "RRGGBB" hex
    # notation used in the OOXML srgbClr val attribute.
    $rr = [Convert]::ToInt32($rrggbb.Substring(0,2), 16)
    $gg = [Convert]::ToInt32($rrggbb.Substring(2,2), 16)
    $bb = [Convert]::ToInt32($rrggbb.Substring(4,2), 16)
    return ($bb * 65536) + ($gg * 256) + $rr
}

$purple = BgrFromRrggbb "7030a0"
$navy   = BgrFromRrggbb "002060"

function FormatRun($rng, $size, $bold, $color) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = $size
    $rng.Font.Bold = $bold
    $rng.Font.Spacing = -0.01
    $rng.Font.Strikethrough = 0
    $rng.Font.Color.RGB = $color
}

function FormatParagraphBase($para) {
    $pf = $para.ParagraphFormat
    $pf.SpaceWithin = 1.0
}

function FormatParagraphBullet($para) {
    $pf = $para.ParagraphFormat
    $pf.SpaceWithin = 1.0
    $pf.SpaceBefore = 561 / 100.0
    $pf.Bullet.Font.Color.RGB = $navy
    $pf.Bullet.Font.Name = "Arial"
    $pf.Bullet.Character = 8226
    $pf.Bullet.Type = 1
}

$p = $ppt.ActivePresentation

# Slide 5 uses the "Blank Slide" layout belonging to the 4th slide master
# (CustomLayouts index 37); the new slide re-uses the very same layout.
$layout = $p.SlideMaster.CustomLayouts.Item(37)
$newSlide = $p.Slides.AddSlide(6, $layout)

# ---------------------------------------------------------------------
# Shape 1 ("CustomShape 1"): slide title
# ---------------------------------------------------------------------
$titleL = EmuToPt 448920
$titleT = EmuToPt 586440
$titleW = EmuToPt 8245440
$titleH = EmuToPt 762840

$shape1 = $newSlide.Shapes.AddShape(1, $titleL, $titleT, $titleW, $titleH)
$shape1.Name = "CustomShape 1"

$tf1 = $shape1.TextFrame
$tf1.MarginLeft = EmuToPt 90000
$tf1.MarginRight = EmuToPt 90000
$tf1.MarginTop = EmuToPt 45000
$tf1.MarginBottom = EmuToPt 45000
$tf1.VerticalAnchor = 3
$tf1.AutoSize = 0

$tf1.TextRange.Text = "Domain Controller"
$titlePara = $tf1.TextRange.Paragraphs(1,1)
FormatParagraphBase $titlePara
FormatRun $titlePara.Characters(1, $titlePara.Length) 36 0 $purple

# ---------------------------------------------------------------------
# Shape 2 ("CustomShape 2"): bulleted body text
# ---------------------------------------------------------------------
$bodyL = EmuToPt 448920
$bodyT = EmuToPt 1502640
$bodyW = EmuToPt 8245440
$bodyH = EmuToPt 3206160

$shape2 = $newSlide.Shapes.AddShape(1, $bodyL, $bodyT, $bodyW, $bodyH)
$shape2.Name = "CustomShape 2"

$tf2 = $shape2.TextFrame
$tf2.MarginLeft = EmuToPt 90000
$tf2.MarginRight = EmuToPt 90000
$tf2.MarginTop = EmuToPt 45000
$tf2.MarginBottom = EmuToPt 45000
$tf2.AutoSize = 2

# Build up the five paragraphs (four bullets + one trailing blank line) as
# plain text first...
$tf2.TextRange.Text = "A server running the AD Domain Service (AD DS) role is called a domain controller."
$tf2.TextRange.InsertAfter([char]13 + "It authenticates and authorizes all users and computers.")
$tf2.TextRange.InsertAfter([char]13 + "A domain controller is contacted when a user logs into a device.")
$tf2.TextRange.InsertAfter([char]13 + "OR accesses another device across the network.")
$tf2.TextRange.InsertAfter([char]13)

# ... then format paragraph 1 (mixed bold runs).
$para1 = $tf2.TextRange.Paragraphs(1,1)
FormatParagraphBullet $para1
$segs = @(
    @{Start=1;  Len=24; Bold=0},
    @{Start=25; Len=14; Bold=-1},
    @{Start=39; Len=2;  Bold=0},
    @{Start=41; Len=5;  Bold=-1},
    @{Start=46; Len=19; Bold=0},
    @{Start=65; Len=17; Bold=-1},
    @{Start=82; Len=1;  Bold=0}
)
foreach ($seg in $segs) {
    FormatRun $para1.Characters($seg.Start, $seg.Len) 28 $seg.Bold $navy
}

# Paragraphs 2-4: single plain run each, all bulleted.
for ($i = 2; $i -le 4; $i++) {
    $para = $tf2.TextRange.Paragraphs($i,1)
    FormatParagraphBullet $para
    FormatRun $para.Characters(1, $para.Length) 28 0 $navy
}

# Paragraph 5: trailing blank paragraph, no bullet.
$para5 = $tf2.TextRange.Paragraphs(5,1)
FormatParagraphBase $para5
$para5.ParagraphFormat.SpaceBefore = 561 / 100.0
